$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.187001824378967
$ws.Range("B1").Value = 2.166004180908203
$ws.Range("C1").Value = 6.333103656768799
$ws.Range("D1").Value = 2.30367374420166
$ws.Range("E1").Value = 1.193690419197083
